$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -1
    3  = 1
    4  = 3
    5  = 3
    6  = 5
    7  = 5
    8  = -3
    10 = -2
    12 = -1
    13 = 2
    14 = 1
    15 = 2
    16 = 3
    17 = 3
    18 = 2
    19 = -1
    20 = -2
    21 = 4
    22 = 2
    23 = -2
    24 = 3
    25 = 3
    26 = 1
    27 = 4
    28 = -4
    29 = -2
    30 = -3
    31 = -3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
